$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("G2").Value = 0.67831
$ws.Range("H2").Value = 1.35662
$ws.Range("I2").Value = 0.004405345910740629
$ws.Range("J2").Value = 0.002944056215323526
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 21.983953
$ws.Range("N2").Value = 65.951859
$ws.Range("O2").Value = 0.3824831516716194
$ws.Range("P2").Value = 0.3824831516716194
$ws.Range("Q2").Value = 14.91193515943
$ws.Range("R2").Value = 89.47161095657999
$ws.Range("S2").Value = 0.001684970588143756
$ws.Range("T2").Value = 0.001126051899935362

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("G3").Value = 0.67831
$ws.Range("H3").Value = 1.35662
$ws.Range("I3").Value = 0.004405345910740629
$ws.Range("J3").Value = 0.002944056215323526
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 35.37535733333333
$ws.Range("N3").Value = 106.126072
$ws.Range("O3").Value = 0.6154706646417534
$ws.Range("P3").Value = 0.6154706646417534
$ws.Range("Q3").Value = 23.99545863277333
$ws.Range("R3").Value = 143.97275179664
$ws.Range("S3").Value = 0.002711361175660365
$ws.Range("T3").Value = 0.001811980235587856

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("G4").Value = 0.67831
$ws.Range("H4").Value = 1.35662
$ws.Range("I4").Value = 0.004405345910740629
$ws.Range("J4").Value = 0.002944056215323526
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.1176083333333333
$ws.Range("N4").Value = 0.352825
$ws.Range("O4").Value = 0.002046183686627228
$ws.Range("P4").Value = 0.002046183686627228
$ws.Range("Q4").Value = 0.07977490858333333
$ws.Range("R4").Value = 0.4786494515
$ws.Range("S4").Value = [double]"9.014146936507443E-06"
$ws.Range("T4").Value = [double]"6.024079800308497E-06"

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 108.0898166666667
$ws.Range("H5").Value = 324.26945
$ws.Range("I5").Value = 0.7019991329115071
$ws.Range("J5").Value = 0.7037103166045329
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 21.983953
$ws.Range("N5").Value = 65.951859
$ws.Range("O5").Value = 0.3824831516716194
$ws.Range("P5").Value = 0.3824831516716194
$ws.Range("Q5").Value = 2376.241449378616
$ws.Range("R5").Value = 21386.17304440755
$ws.Range("S5").Value = 0.2685028408267373
$ws.Range("T5").Value = 0.2691573397587349

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 108.0898166666667
$ws.Range("H6").Value = 324.26945
$ws.Range("I6").Value = 0.7019991329115071
$ws.Range("J6").Value = 0.7037103166045329
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 35.37535733333333
$ws.Range("N6").Value = 106.126072
$ws.Range("O6").Value = 0.6154706646417534
$ws.Range("P6").Value = 0.6154706646417534
$ws.Range("Q6").Value = 3823.715888677822
$ws.Range("R6").Value = 34413.4429981004
$ws.Range("S6").Value = 0.4320598729109798
$ws.Range("T6").Value = 0.4331130562758506

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 108.0898166666667
$ws.Range("H7").Value = 324.26945
$ws.Range("I7").Value = 0.7019991329115071
$ws.Range("J7").Value = 0.7037103166045329
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.1176083333333333
$ws.Range("N7").Value = 0.352825
$ws.Range("O7").Value = 0.002046183686627228
$ws.Range("P7").Value = 0.002046183686627228
$ws.Range("Q7").Value = 12.71226318847222
$ws.Range("R7").Value = 114.41036869625
$ws.Range("S7").Value = 0.001436419173789985
$ws.Range("T7").Value = 0.001439920569947477

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 17.58332566666667
$ws.Range("H8").Value = 52.749977
$ws.Range("I8").Value = 0.114196505761187
$ws.Range("J8").Value = 0.1144748696355818
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 21.983953
$ws.Range("N8").Value = 65.951859
$ws.Range("O8").Value = 0.3824831516716194
$ws.Range("P8").Value = 0.3824831516716194
$ws.Range("Q8").Value = 386.5510050396937
$ws.Range("R8").Value = 3478.959045357243
$ws.Range("S8").Value = 0.04367823943342505
$ws.Range("T8").Value = 0.04378470892541511

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 17.58332566666667
$ws.Range("H9").Value = 52.749977
$ws.Range("I9").Value = 0.114196505761187
$ws.Range("J9").Value = 0.1144748696355818
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 35.37535733333333
$ws.Range("N9").Value = 106.126072
$ws.Range("O9").Value = 0.6154706646417534
$ws.Range("P9").Value = 0.6154706646417534
$ws.Range("Q9").Value = 622.0164285667049
$ws.Range("R9").Value = 5598.147857100344
$ws.Range("S9").Value = 0.07028459930060359
$ws.Range("T9").Value = 0.07045592409938962

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 17.58332566666667
$ws.Range("H10").Value = 52.749977
$ws.Range("I10").Value = 0.114196505761187
$ws.Range("J10").Value = 0.1144748696355818
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 0.1176083333333333
$ws.Range("N10").Value = 0.352825
$ws.Range("O10").Value = 0.002046183686627228
$ws.Range("P10").Value = 0.002046183686627228
$ws.Range("Q10").Value = 2.067945626113889
$ws.Range("R10").Value = 18.611510635025
$ws.Range("S10").Value = 0.0002336670271583731
$ws.Range("T10").Value = 0.0002342366107771062

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 25.99611466666667
$ws.Range("H11").Value = 77.988344
$ws.Range("I11").Value = 0.1688341281153816
$ws.Range("J11").Value = 0.1692456759269281
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 21.983953
$ws.Range("N11").Value = 65.951859
$ws.Range("O11").Value = 0.3824831516716194
$ws.Range("P11").Value = 0.3824831516716194
$ws.Range("Q11").Value = 571.4973630146106
$ws.Range("R11").Value = 5143.476267131496
$ws.Range("S11").Value = 0.06457620943130113
$ws.Range("T11").Value = 0.064733619535325

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 25.99611466666667
$ws.Range("H12").Value = 77.988344
$ws.Range("I12").Value = 0.1688341281153816
$ws.Range("J12").Value = 0.1692456759269281
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 35.37535733333333
$ws.Range("N12").Value = 106.126072
$ws.Range("O12").Value = 0.6154706646417534
$ws.Range("P12").Value = 0.6154706646417534
$ws.Range("Q12").Value = 919.6218456116409
$ws.Range("R12").Value = 8276.596610504766
$ws.Range("S12").Value = 0.1039124530453849
$ws.Range("T12").Value = 0.1041657486504892

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 25.99611466666667
$ws.Range("H13").Value = 77.988344
$ws.Range("I13").Value = 0.1688341281153816
$ws.Range("J13").Value = 0.1692456759269281
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.1176083333333333
$ws.Range("N13").Value = 0.352825
$ws.Range("O13").Value = 0.002046183686627228
$ws.Range("P13").Value = 0.002046183686627228
$ws.Range("Q13").Value = 3.057359719088889
$ws.Range("R13").Value = 27.5162374718
$ws.Range("S13").Value = 0.0003454656386956252
$ws.Range("T13").Value = 0.0003463077411138789

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 1.181792
$ws.Range("H14").Value = 3.545376
$ws.Range("I14").Value = 0.007675255494605696
$ws.Range("J14").Value = 0.00769396459469775
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 21.983953
$ws.Range("N14").Value = 65.951859
$ws.Range("O14").Value = 0.3824831516716194
$ws.Range("P14").Value = 0.3824831516716194
$ws.Range("Q14").Value = 25.980459783776
$ws.Range("R14").Value = 233.824138053984
$ws.Range("S14").Value = 0.002935655911461701
$ws.Range("T14").Value = 0.00294281182702985

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 1.181792
$ws.Range("H15").Value = 3.545376
$ws.Range("I15").Value = 0.007675255494605696
$ws.Range("J15").Value = 0.00769396459469775
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 35.37535733333333
$ws.Range("N15").Value = 106.126072
$ws.Range("O15").Value = 0.6154706646417534
$ws.Range("P15").Value = 0.6154706646417534
$ws.Range("Q15").Value = 41.80631429367467
$ws.Range("R15").Value = 376.256828643072
$ws.Range("S15").Value = 0.004723894600560237
$ws.Range("T15").Value = 0.004735409502828743

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 1.181792
$ws.Range("H16").Value = 3.545376
$ws.Range("I16").Value = 0.007675255494605696
$ws.Range("J16").Value = 0.00769396459469775
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.1176083333333333
$ws.Range("N16").Value = 0.352825
$ws.Range("O16").Value = 0.002046183686627228
$ws.Range("P16").Value = 0.002046183686627228
$ws.Range("Q16").Value = 0.1389885874666666
$ws.Range("R16").Value = 1.2508972872
$ws.Range("S16").Value = [double]"1.570498258375817E-05"
$ws.Range("T16").Value = [double]"1.574326483915801E-05"

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("G17").Value = 0.444929
$ws.Range("H17").Value = 0.889858
$ws.Range("I17").Value = 0.002889631806577992
$ws.Range("J17").Value = 0.001931117022935945
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 21.983953
$ws.Range("N17").Value = 65.951859
$ws.Range("O17").Value = 0.3824831516716194
$ws.Range("P17").Value = 0.3824831516716194
$ws.Range("Q17").Value = 9.781298224337
$ws.Range("R17").Value = 58.687789346022
$ws.Range("S17").Value = 0.001105235480550506
$ws.Range("T17").Value = 0.0007386197251792554

# Row 18
$ws.Range("E18").Value = 2
$ws.Range("G18").Value = 0.444929
$ws.Range("H18").Value = 0.889858
$ws.Range("I18").Value = 0.002889631806577992
$ws.Range("J18").Value = 0.001931117022935945
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 35.37535733333333
$ws.Range("N18").Value = 106.126072
$ws.Range("O18").Value = 0.6154706646417534
$ws.Range("P18").Value = 0.6154706646417534
$ws.Range("Q18").Value = 15.73952236296267
$ws.Range("R18").Value = 94.437134177776
$ws.Range("S18").Value = 0.001778483608564507
$ws.Range("T18").Value = 0.00118854587760739

# Row 19
$ws.Range("E19").Value = 2
$ws.Range("G19").Value = 0.444929
$ws.Range("H19").Value = 0.889858
$ws.Range("I19").Value = 0.002889631806577992
$ws.Range("J19").Value = 0.001931117022935945
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 0.1176083333333333
$ws.Range("N19").Value = 0.352825
$ws.Range("O19").Value = 0.002046183686627228
$ws.Range("P19").Value = 0.002046183686627228
$ws.Range("Q19").Value = 0.05232735814166667
$ws.Range("R19").Value = 0.31396414885
$ws.Range("S19").Value = [double]"5.912717462979052E-06"
$ws.Range("T19").Value = [double]"3.95142014929967E-06"
